$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the date serial 45207 (2023-10-08) for every
# data row (2 through 83). The update bumps it by one day to 45208 (2023-10-09).
for ($r = 2; $r -le 83; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value2 = 45208
    }
}
